$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "328.42"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.01%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "44.16"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.54%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.487"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-1.15%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08040"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.04%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.032"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "7.15%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9547"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.12%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.1105"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-7.01%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1875"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.63%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "10.15"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09953"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "2.57%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04746"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "5.83%"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.66%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001274"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.82%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.04083"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-2.47%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005873"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.59%"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.82%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.416"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "3.62%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "4.75%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3411"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.46%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1401"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.97%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "5.42%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004339"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.70%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.47%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0003747"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-6.13%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02573"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "-2.71%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05673"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "3.13%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007708"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.93%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1399"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.39%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007349"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-10.21%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "0.64%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008509"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-3.35%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00007119"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "0.13%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.11%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0005811"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.01%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003504"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "54.28%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.003506"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "3.81%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002104"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.11%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002004"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.11%"
